# Add two new columns, I ("I0") and J ("IF"), to the sheet.
# Header row (row 1) gets the new labels with the same style (s="1") as
# the other header cells (bold, centered, bordered).
# Data rows 2-70 get the corresponding numeric values from the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cells -----------------------------------------------------
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
# Match the formatting used by the other header cells (B1:H1) - bold,
# centered, bordered - by copying the format from an existing header cell.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Data values --------------------------------------------------------
$iValues = @(8,6,7,6,7,9,8,10,7,5,6,8,10,8,6,9,9,9,7,6,7,7,9,8,5,9,8,9,7,8,8,10,9,7,8,6,9,9,6,9,6,8,8,8,9,9,7,7,9,8,8,5,6,5,8,4,8,9,6,6,9,5,5,5,6,6,7,3,4)
$jValues = @(8,7,7,7,7,9,8,10,7,6,6,8,10,8,6,9,9,9,7,6,8,7,9,8,6,9,8,9,7,8,9,10,9,7,9,7,9,9,6,9,7,8,8,8,9,9,7,8,9,8,8,6,7,6,9,5,8,9,6,6,10,6,6,6,7,6,8,3,4)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
